$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.169.23'
$ws.Range("E2").Value = '  +3.07%  '
$ws.Range("D3").Value = '2.363.86'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '''310.68'
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("D6").Value = '''108.17'
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").Value = '''0.631'
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("D9").Value = '''0.606'
$ws.Range("E9").Value = '  -1.99%  '
$ws.Range("D10").Value = '''40.80'
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("D12").Value = '''8.43'
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("E13").Value = '  +1.30%  '
$ws.Range("D14").Value = '''0.976'
$ws.Range("E14").Value = '  -2.90%  '
$ws.Range("D15").Value = '2.719.46'
$ws.Range("E15").Value = '  +1.00%  '
$ws.Range("D16").Value = '''15.18'
$ws.Range("E16").Value = '  -1.71%  '
$ws.Range("D17").Value = '2.358.48'
$ws.Range("E17").Value = '  +1.09%  '
$ws.Range("D18").Value = '45.059.79'
$ws.Range("E18").Value = '  +3.02%  '
$ws.Range("D19").Value = '''14.49'
$ws.Range("E19").Value = '  +10.90%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '''7.19'
$ws.Range("E20").Value = '  -4.65%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '''0.0000106'
$ws.Range("E21").Value = '  -0.94%  '
$ws.Range("D22").Value = '''72.95'
$ws.Range("E22").Value = '  -1.66%  '
$ws.Range("D23").Value = '''3.51'
$ws.Range("E23").Value = '  +1.09%  '
$ws.Range("D24").Value = '''258.90'
$ws.Range("E24").Value = '  -3.63%  '
$ws.Range("E25").Value = '  +1.51%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("E27").Value = '  -0.45%  '
$ws.Range("D28").Value = '''7.18'
$ws.Range("E28").Value = '  -6.57%  '
$ws.Range("E29").Value = '  +1.16%  '
$ws.Range("E30").Value = '  +9.35%  '
$ws.Range("D31").Value = '''22.29'
$ws.Range("E31").Value = '  -1.35%  '
$ws.Range("D32").Value = '''37.11'
$ws.Range("E32").Value = '  -5.02%  '
$ws.Range("D33").Value = '''167.87'
$ws.Range("E33").Value = '  -0.40%  '
$ws.Range("E34").Value = '  +5.46%  '
$ws.Range("E35").Value = '  -1.34%  '
$ws.Range("E36").Value = '  +1.40%  '
$ws.Range("D37").Value = '''4.67'
$ws.Range("E37").Value = '  -0.93%  '
$ws.Range("D38").Value = '''3.96'
$ws.Range("E38").Value = '  +4.24%  '
$ws.Range("D39").Value = '''2.90'
$ws.Range("E39").Value = '  +0.56%  '
$ws.Range("D40").Value = '''0.0351'
$ws.Range("E40").Value = '  -3.16%  '
$ws.Range("D41").Value = '''1.76'
$ws.Range("E41").Value = '  +3.08%  '
$ws.Range("D42").Value = '''100.09'
$ws.Range("E42").Value = '  -4.74%  '
$ws.Range("D43").Value = '1.896.79'
$ws.Range("E43").Value = '  +13.58%  '
$ws.Range("D44").Value = '''69.28'
$ws.Range("E44").Value = '  -3.35%  '
$ws.Range("D45").Value = '''0.229'
$ws.Range("E45").Value = '  -4.26%  '
$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").Value = '''12.83'
$ws.Range("E46").Value = '  -4.63%  '
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").Value = '''1.00'
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("D48").Value = '''81.75'
$ws.Range("E48").Value = '  +6.52%  '
$ws.Range("E49").Value = '  +8.04%  '
$ws.Range("D50").Value = '''110.28'
$ws.Range("E50").Value = '  -3.23%  '
$ws.Range("D51").Value = '''9.15'
$ws.Range("E51").Value = '  +2.41%  '
